$d = $word.ActiveDocument

# Locate the target paragraph robustly via a unique substring instead of a
# hard-coded paragraph index.
$finder = $d.Content
$ok = $finder.Find.Execute("constelação de Perseu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $finder.Paragraphs(1)
$full = $para.Range

# New, single-run plain-text replacement for the whole paragraph body
# (matches the target: all the per-word runs/proofErr wrappers/fonts are
# collapsed into one plain run).
$newText = "Está a participar numa campanha global para observar e registar as estrelas mais fracas visíveis como forma de medir a poluição luminosa num determinado local. Localizando e observando a  constelação de leão no céu noturno e,  comparando-a com cartas estelares, pessoas de todo o mundo aprenderão  como as luzes da sua comunidade contribuem para a poluição luminosa. As suas contribuições para a base de dados on-line irão documentar a visibilidade do céu noturno em todo o mundo."

# Replace the paragraph's text (excluding its trailing paragraph mark) with
# the new plain text, preserving the paragraph element itself (its pPr,
# paraId, rsidP, ...).
$body = $d.Range($full.Start, $full.End - 1)
$body.Delete()
$body.InsertAfter($newText)
